# Edit commit: Tue, Jul 28, 2020 3:05:54 AM
#
# 1) Re-apply the built-in "Medium Style 2 - Accent 1" table style
#    ({C3E50673-3C97-4C2B-88B2-FEB6AF250356}) to the three data tables
#    (previously the generic Google-Slides-imported "Table_0" style,
#    {DAA98444-3182-41E9-8479-C54D316EB2F9}).
# 2) Swap the deck's colour theme from "Integral" (Red Violet) to the
#    stock "Office Theme" (Office) palette.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$oldTableStyle = "{DAA98444-3182-41E9-8479-C54D316EB2F9}"
$newTableStyle = "{C3E50673-3C97-4C2B-88B2-FEB6AF250356}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable -and $shape.Table.Style -eq $oldTableStyle) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# index : role      : new "Office" RGB
#   1   : dk1       : 000000
#   2   : lt1       : FFFFFF
#   3   : dk2       : 44546A
#   4   : lt2       : E7E6E6
#   5   : accent1   : 5B9BD5
#   6   : accent2   : ED7D31
#   7   : accent3   : A5A5A5
#   8   : accent4   : FFC000
#   9   : accent5   : 4472C4
#  10   : accent6   : 70AD47
#  11   : hlink     : 0563C1
#  12   : folHlink  : 954F72
$officeThemeColors = @{
    1  = 0
    2  = 16777215
    3  = 6968388
    4  = 15132391
    5  = 13998939
    6  = 3243501
    7  = 10855845
    8  = 49407
    9  = 12874308
    10 = 4697456
    11 = 12673797
    12 = 7491477
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i]
}
